# 自动更新Excel文件 - 2026-02-04 23:19:31
#
# Daily "refill countdown" refresh.
#   D = total cycle length in days (总天)
#   E = days remaining (剩余)
#   F = cycle start date, stored as an integer literal yyyyMMdd (开始时间)
#
# For every data row, the number of days elapsed since F is (D - E); that
# elapsed count is anchored to a single "as-of" date shared by the whole
# sheet. Rolling the as-of date forward by one day means every row's
# remaining-days count drops by one EXCEPT rows that had already bottomed
# out at E = 1 the day before: for those, a fresh cycle restarts on the new
# as-of date (F resets to the new as-of date, E resets to the full D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

function Parse-Ymd([double]$n) {
    $y = [Math]::Floor($n / 10000)
    $m = [Math]::Floor(($n - $y * 10000) / 100)
    $d = $n - $y * 10000 - $m * 100

    if ($y -lt 1900 -or $y -gt 9999) { return $null }
    if ($m -lt 1 -or $m -gt 12) { return $null }
    if ($d -lt 1 -or $d -gt 31) { return $null }

    $dt = Get-Date -Year ([int]$y) -Month ([int]$m) -Day ([int]$d) -Hour 0 -Minute 0 -Second 0 -Millisecond 0
    if ($dt.Year -ne $y -or $dt.Month -ne $m -or $dt.Day -ne $d) { return $null }
    return $dt
}

for ($row = 2; $row -le $lastRow; $row++) {
    $dCell = $ws.Cells.Item($row, 4)   # D: 总天
    $eCell = $ws.Cells.Item($row, 5)   # E: 剩余
    $fCell = $ws.Cells.Item($row, 6)   # F: 开始时间 (yyyyMMdd)

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($dVal -eq $null -or $eVal -eq $null -or $fVal -eq $null) {
        continue
    }

    $startDate = Parse-Ymd([double]$fVal)
    if ($startDate -eq $null) {
        # Malformed start date (e.g. a fat-fingered yyyyMMdd) - leave as-is,
        # matches the one row in this sheet that was left untouched.
        continue
    }

    $total = [int]$dVal
    $remainingOld = [int]$eVal
    $elapsedOld = $total - $remainingOld
    $asOfOld = $startDate.AddDays($elapsedOld)
    $asOfNew = $asOfOld.AddDays(1)

    if ($remainingOld -eq 1) {
        # Cycle had bottomed out - it restarts today.
        $newStart = $asOfNew
        $remainingNew = $total
    } else {
        $newStart = $startDate
        $elapsedNew = [int]($asOfNew.ToOADate() - $newStart.ToOADate())
        $remainingNew = $total - $elapsedNew
    }

    $newYmd = $newStart.Year * 10000 + $newStart.Month * 100 + $newStart.Day

    $eCell.Value = $remainingNew
    $fCell.Value = $newYmd
}
